# Applies the cryptos.xlsx data refresh described in the commit message
# "Updated cryptos list on Wed Sep 11 16:55:13 UTC 2024 with GitHub Actions".
# Price/volume figures are refreshed per-row; two rows (21/22) also swap
# which coin (Polkadot vs Uniswap) occupies that rank.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.044.13"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").Value = "2.319.31"
$ws.Range("E3").Value = "  -0.53%  "
$ws.Range("D4").Formula = "=""0.998"""
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D4").PasteSpecial(-4163) | Out-Null
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Formula = "=""522.64"""
$ws.Range("D5").Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4163) | Out-Null
$ws.Range("D6").Formula = "=""131.94"""
$ws.Range("D6").Copy() | Out-Null
$ws.Range("D6").PasteSpecial(-4163) | Out-Null
$ws.Range("E6").Value = "  -2.07%  "
$ws.Range("D7").Formula = "=""0.996"""
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D7").PasteSpecial(-4163) | Out-Null
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("E8").Value = "  -0.66%  "
$ws.Range("D9").Value = "2.338.93"
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("E10").Value = "  -1.54%  "
$ws.Range("D12").Formula = "=""5.30"""
$ws.Range("D12").Copy() | Out-Null
$ws.Range("D12").PasteSpecial(-4163) | Out-Null
$ws.Range("E12").Value = "  -1.28%  "
$ws.Range("E13").Value = "  -0.16%  "
$ws.Range("D14").Formula = "=""23.51"""
$ws.Range("D14").Copy() | Out-Null
$ws.Range("D14").PasteSpecial(-4163) | Out-Null
$ws.Range("E14").Value = "  -1.63%  "
$ws.Range("D15").Value = "2.736.07"
$ws.Range("E15").Value = "  -0.43%  "
$ws.Range("D16").Value = "56.880.49"
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("E17").Value = "  -1.16%  "
$ws.Range("D18").Value = "2.339.97"
$ws.Range("E18").Value = "  +0.37%  "
$ws.Range("D19").Formula = "=""336.04"""
$ws.Range("D19").Copy() | Out-Null
$ws.Range("D19").PasteSpecial(-4163) | Out-Null
$ws.Range("E19").Value = "  +2.98%  "
$ws.Range("E20").Value = "  -0.78%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Formula = "=""6.89"""
$ws.Range("D21").Copy() | Out-Null
$ws.Range("D21").PasteSpecial(-4163) | Out-Null
$ws.Range("E21").Value = "  +4.20%  "
$ws.Range("B22").Value = "Polkadot"
$ws.Range("C22").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D22").Formula = "=""4.15"""
$ws.Range("D22").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4163) | Out-Null
$ws.Range("E22").Value = "  -1.35%  "
$ws.Range("D23").Formula = "=""0.999"""
$ws.Range("D23").Copy() | Out-Null
$ws.Range("D23").PasteSpecial(-4163) | Out-Null
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").Formula = "=""61.59"""
$ws.Range("D24").Copy() | Out-Null
$ws.Range("D24").PasteSpecial(-4163) | Out-Null
$ws.Range("E24").Value = "  +1.23%  "
$ws.Range("E25").Value = "  +9.23%  "
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("D27").Formula = "=""0.998"""
$ws.Range("D27").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4163) | Out-Null
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("E28").Value = "  +1.96%  "
$ws.Range("E29").Value = "  -0.40%  "
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("E31").Value = "  -2.17%  "
$ws.Range("D32").Formula = "=""6.10"""
$ws.Range("D32").Copy() | Out-Null
$ws.Range("D32").PasteSpecial(-4163) | Out-Null
$ws.Range("E32").Value = "  -1.60%  "
$ws.Range("E33").Value = "  -0.43%  "
$ws.Range("D35").Formula = "=""0.996"""
$ws.Range("D35").Copy() | Out-Null
$ws.Range("D35").PasteSpecial(-4163) | Out-Null
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").Formula = "=""1.26"""
$ws.Range("D36").Copy() | Out-Null
$ws.Range("D36").PasteSpecial(-4163) | Out-Null
$ws.Range("E36").Value = "  -0.52%  "
$ws.Range("D37").Formula = "=""3.99"""
$ws.Range("D37").Copy() | Out-Null
$ws.Range("D37").PasteSpecial(-4163) | Out-Null
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("D38").Formula = "=""0.899"""
$ws.Range("D38").Copy() | Out-Null
$ws.Range("D38").PasteSpecial(-4163) | Out-Null
$ws.Range("E38").Value = "  -1.49%  "
$ws.Range("E39").Value = "  +1.09%  "
$ws.Range("D40").Formula = "=""38.87"""
$ws.Range("D40").Copy() | Out-Null
$ws.Range("D40").PasteSpecial(-4163) | Out-Null
$ws.Range("E40").Value = "  +1.41%  "
$ws.Range("D41").Formula = "=""149.05"""
$ws.Range("D41").Copy() | Out-Null
$ws.Range("D41").PasteSpecial(-4163) | Out-Null
$ws.Range("E41").Value = "  +4.08%  "
$ws.Range("D42").Formula = "=""0.374"""
$ws.Range("D42").Copy() | Out-Null
$ws.Range("D42").PasteSpecial(-4163) | Out-Null
$ws.Range("E42").Value = "  -1.42%  "
$ws.Range("D43").Formula = "=""286.38"""
$ws.Range("D43").Copy() | Out-Null
$ws.Range("D43").PasteSpecial(-4163) | Out-Null
$ws.Range("E43").Value = "  +3.31%  "
$ws.Range("E44").Value = "  -0.72%  "
$ws.Range("D45").Formula = "=""5.08"""
$ws.Range("D45").Copy() | Out-Null
$ws.Range("D45").PasteSpecial(-4163) | Out-Null
$ws.Range("E45").Value = "  -1.95%  "
$ws.Range("E46").Value = "  -0.67%  "
$ws.Range("E47").Value = "  -0.87%  "
$ws.Range("D48").Formula = "=""0.558"""
$ws.Range("D48").Copy() | Out-Null
$ws.Range("D48").PasteSpecial(-4163) | Out-Null
$ws.Range("E48").Value = "  -0.64%  "
$ws.Range("D49").Formula = "=""18.48"""
$ws.Range("D49").Copy() | Out-Null
$ws.Range("D49").PasteSpecial(-4163) | Out-Null
$ws.Range("E49").Value = "  +2.75%  "
$ws.Range("E50").Value = "  -1.35%  "
$ws.Range("D51").Formula = "=""0.376"""
$ws.Range("D51").Copy() | Out-Null
$ws.Range("D51").PasteSpecial(-4163) | Out-Null
$ws.Range("E51").Value = "  -0.78%  "

$excel.CutCopyMode = 0

